# Generate Report for Handoff
#
# Rotates the localization-status report to a freshly generated source file
# (2470d0a0-52c6-4b88-909c-edf929b018ad.md replacing
# d86cc8b9-4bff-4efe-9b6f-dd23c4759863.md), refreshes the "Latest Handoff"
# file/datetime stamps for both locales, and clears out the stale
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# values now that a new handoff cycle has started.

$wb = $excel.ActiveWorkbook

$oldId = "d86cc8b9-4bff-4efe-9b6f-dd23c4759863"
$newId = "2470d0a0-52c6-4b88-909c-edf929b018ad"
$oldHash = "23e881281bce6773eecbbf5ac61505876f7789d3"
$newHash = "a34d8a292b3aab8e612d7fe91cf63f4b70ca5aca"

$newMdName = $newId + ".md"
$newMdPath = "e2e\" + $newMdName

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMdPath
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = $newId + "." + $newHash + ".zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-03 21:02:36"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Column -eq 9) {
        $hl.Delete()
    } else {
        $hl.TextToDisplay = $newMdName
    }
}

$wsZh.Columns.Item(9).ColumnWidth = 17.83
$wsZh.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = $newId + "." + $newHash + ".de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-03 21:02:40"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Column -eq 9) {
        $hl.Delete()
    } else {
        $hl.TextToDisplay = $newMdName
    }
}

$wsDe.Columns.Item(9).ColumnWidth = 17.83
$wsDe.Columns.Item(10).ColumnWidth = 20.83
